# Update column G ("K" - strikeouts) values on Sheet1 to reflect the
# regenerated save_data (switched from "Strike#" to "K", recalculated).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 3
    3  = 5
    4  = 1
    5  = 4
    6  = 5
    7  = 3
    8  = 8
    9  = 5
    10 = 9
    11 = 8
    12 = 7
    13 = 8
    14 = 5
    15 = 6
    16 = 10
    17 = 3
    18 = 12
    19 = 5
    20 = 4
    21 = 8
    22 = 4
    23 = 3
    24 = 4
    25 = 4
    26 = 6
    27 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}

$wb.Save()
